$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 42.05115733333333
$ws.Range("H2").Value = 126.153472
$ws.Range("I2").Value = 0.1594435451835853
$ws.Range("J2").Value = 0.1594435451835853
$ws.Range("M2").Value = 1.175645333333333
$ws.Range("N2").Value = 3.526936
$ws.Range("O2").Value = 0.06804514706690673
$ws.Range("P2").Value = 0.06804514706690673
$ws.Range("Q2").Value = 49.43724688019911
$ws.Range("R2").Value = 444.935221921792
$ws.Range("S2").Value = 0.01084935948088605
$ws.Range("T2").Value = 0.01084935948088605
$ws.Range("G3").Value = 42.05115733333333
$ws.Range("H3").Value = 126.153472
$ws.Range("I3").Value = 0.1594435451835853
$ws.Range("J3").Value = 0.1594435451835853
$ws.Range("O3").Value = 0.5504564499973018
$ws.Range("P3").Value = 0.5504564499973019
$ws.Range("Q3").Value = 399.926410454471
$ws.Range("R3").Value = 3599.337694090239
$ws.Range("S3").Value = 0.08776672785674074
$ws.Range("T3").Value = 0.08776672785674078
$ws.Range("G4").Value = 42.05115733333333
$ws.Range("H4").Value = 126.153472
$ws.Range("I4").Value = 0.1594435451835853
$ws.Range("J4").Value = 0.1594435451835853
$ws.Range("M4").Value = 2.410127666666666
$ws.Range("N4").Value = 7.230383
$ws.Range("O4").Value = 0.1394957193964002
$ws.Range("P4").Value = 0.1394957193964002
$ws.Range("Q4").Value = 101.3486577044195
$ws.Range("R4").Value = 912.137919339776
$ws.Range("S4").Value = 0.02224169203849667
$ws.Range("T4").Value = 0.02224169203849668
$ws.Range("G5").Value = 42.05115733333333
$ws.Range("H5").Value = 126.153472
$ws.Range("I5").Value = 0.1594435451835853
$ws.Range("J5").Value = 0.1594435451835853
$ws.Range("M5").Value = 4.181184666666667
$ws.Range("N5").Value = 12.543554
$ws.Range("O5").Value = 0.2420026835393912
$ws.Range("P5").Value = 0.2420026835393912
$ws.Range("Q5").Value = 175.8236542577209
$ws.Range("R5").Value = 1582.412888319488
$ws.Range("S5").Value = 0.03858576580746181
$ws.Range("T5").Value = 0.03858576580746181
$ws.Range("G6").Value = 57.66057933333332
$ws.Range("I6").Value = 0.2186291119973147
$ws.Range("J6").Value = 0.2186291119973148
$ws.Range("M6").Value = 1.175645333333333
$ws.Range("N6").Value = 3.526936
$ws.Range("O6").Value = 0.06804514706690673
$ws.Range("P6").Value = 0.06804514706690673
$ws.Range("Q6").Value = 67.78839101052978
$ws.Range("R6").Value = 610.095519094768
$ws.Range("S6").Value = 0.01487665007896451
$ws.Range("T6").Value = 0.01487665007896451
$ws.Range("G7").Value = 57.66057933333332
$ws.Range("I7").Value = 0.2186291119973147
$ws.Range("J7").Value = 0.2186291119973148
$ws.Range("O7").Value = 0.5504564499973018
$ws.Range("P7").Value = 0.5504564499973019
$ws.Range("Q7").Value = 548.3794021342176
$ws.Range("R7").Value = 4935.414619207959
$ws.Range("S7").Value = 0.1203458048561044
$ws.Range("T7").Value = 0.1203458048561044
$ws.Range("G8").Value = 57.66057933333332
$ws.Range("I8").Value = 0.2186291119973147
$ws.Range("J8").Value = 0.2186291119973148
$ws.Range("M8").Value = 2.410127666666666
$ws.Range("N8").Value = 7.230383
$ws.Range("O8").Value = 0.1394957193964002
$ws.Range("P8").Value = 0.1394957193964002
$ws.Range("Q8").Value = 138.9693575272948
$ws.Range("R8").Value = 1250.724217745654
$ws.Range("S8").Value = 0.03049782525906157
$ws.Range("T8").Value = 0.03049782525906158
$ws.Range("G9").Value = 57.66057933333332
$ws.Range("I9").Value = 0.2186291119973147
$ws.Range("J9").Value = 0.2186291119973148
$ws.Range("M9").Value = 4.181184666666667
$ws.Range("N9").Value = 12.543554
$ws.Range("O9").Value = 0.2420026835393912
$ws.Range("P9").Value = 0.2420026835393912
$ws.Range("Q9").Value = 241.0895301796502
$ws.Range("R9").Value = 2169.805771616852
$ws.Range("S9").Value = 0.05290883180318427
$ws.Range("T9").Value = 0.05290883180318428
$ws.Range("G10").Value = 99.15200299999999
$ws.Range("H10").Value = 297.456009
$ws.Range("I10").Value = 0.3759503393701321
$ws.Range("J10").Value = 0.3759503393701321
$ws.Range("M10").Value = 1.175645333333333
$ws.Range("N10").Value = 3.526936
$ws.Range("O10").Value = 0.06804514706690673
$ws.Range("P10").Value = 0.06804514706690673
$ws.Range("Q10").Value = 116.5675896176027
$ws.Range("R10").Value = 1049.108306558424
$ws.Range("S10").Value = 0.02558159613229413
$ws.Range("T10").Value = 0.02558159613229414
$ws.Range("G11").Value = 99.15200299999999
$ws.Range("H11").Value = 297.456009
$ws.Range("I11").Value = 0.3759503393701321
$ws.Range("J11").Value = 0.3759503393701321
$ws.Range("O11").Value = 0.5504564499973018
$ws.Range("P11").Value = 0.5504564499973019
$ws.Range("Q11").Value = 942.9824804780865
$ws.Range("R11").Value = 8486.842324302779
$ws.Range("S11").Value = 0.2069442891849637
$ws.Range("T11").Value = 0.2069442891849638
$ws.Range("G12").Value = 99.15200299999999
$ws.Range("H12").Value = 297.456009
$ws.Range("I12").Value = 0.3759503393701321
$ws.Range("J12").Value = 0.3759503393701321
$ws.Range("M12").Value = 2.410127666666666
$ws.Range("N12").Value = 7.230383
$ws.Range("O12").Value = 0.1394957193964002
$ws.Range("P12").Value = 0.1394957193964002
$ws.Range("Q12").Value = 238.9689856357163
$ws.Range("R12").Value = 2150.720870721447
$ws.Range("S12").Value = 0.05244346304775737
$ws.Range("T12").Value = 0.05244346304775739
$ws.Range("G13").Value = 99.15200299999999
$ws.Range("H13").Value = 297.456009
$ws.Range("I13").Value = 0.3759503393701321
$ws.Range("J13").Value = 0.3759503393701321
$ws.Range("M13").Value = 4.181184666666667
$ws.Range("N13").Value = 12.543554
$ws.Range("O13").Value = 0.2420026835393912
$ws.Range("P13").Value = 0.2420026835393912
$ws.Range("Q13").Value = 414.5728346128873
$ws.Range("R13").Value = 3731.155511515986
$ws.Range("S13").Value = 0.09098099100511678
$ws.Range("T13").Value = 0.09098099100511681
$ws.Range("G14").Value = 64.87322933333333
$ws.Range("H14").Value = 194.619688
$ws.Range("I14").Value = 0.2459770034489679
$ws.Range("J14").Value = 0.2459770034489679
$ws.Range("M14").Value = 1.175645333333333
$ws.Range("N14").Value = 3.526936
$ws.Range("O14").Value = 0.06804514706690673
$ws.Range("P14").Value = 0.06804514706690673
$ws.Range("Q14").Value = 76.26790932399645
$ws.Range("R14").Value = 686.411183915968
$ws.Range("S14").Value = 0.01673754137476204
$ws.Range("T14").Value = 0.01673754137476204
$ws.Range("G15").Value = 64.87322933333333
$ws.Range("H15").Value = 194.619688
$ws.Range("I15").Value = 0.2459770034489679
$ws.Range("J15").Value = 0.2459770034489679
$ws.Range("O15").Value = 0.5504564499973018
$ws.Range("P15").Value = 0.5504564499973019
$ws.Range("Q15").Value = 616.9751176218842
$ws.Range("R15").Value = 5552.776058596959
$ws.Range("S15").Value = 0.1353996280994929
$ws.Range("T15").Value = 0.135399628099493
$ws.Range("G16").Value = 64.87322933333333
$ws.Range("H16").Value = 194.619688
$ws.Range("I16").Value = 0.2459770034489679
$ws.Range("J16").Value = 0.2459770034489679
$ws.Range("M16").Value = 2.410127666666666
$ws.Range("N16").Value = 7.230383
$ws.Range("O16").Value = 0.1394957193964002
$ws.Range("P16").Value = 0.1394957193964002
$ws.Range("Q16").Value = 156.3527648422782
$ws.Range("R16").Value = 1407.174883580504
$ws.Range("S16").Value = 0.03431273905108459
$ws.Range("T16").Value = 0.0343127390510846
$ws.Range("G17").Value = 64.87322933333333
$ws.Range("H17").Value = 194.619688
$ws.Range("I17").Value = 0.2459770034489679
$ws.Range("J17").Value = 0.2459770034489679
$ws.Range("M17").Value = 4.181184666666667
$ws.Range("N17").Value = 12.543554
$ws.Range("O17").Value = 0.2420026835393912
$ws.Range("P17").Value = 0.2420026835393912
$ws.Range("Q17").Value = 271.2469517656835
$ws.Range("R17").Value = 2441.222565891152
$ws.Range("S17").Value = 0.0595270949236283
$ws.Range("T17").Value = 0.05952709492362832
